# Update the weekly menu: shift dates by one week and refresh dish names.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date range for the menu (moves from 2025-06-02..06-06 to 2025-06-09..06-13)
$ws.Range("A2").Value = 45817
$ws.Range("B2").Value = 45821

# Main dishes (column C)
$ws.Range("C2").Value = "Vepřová panenka s pepřovou omáčkou a knedlíkem / Traditional Hungarian goulash with pork meat and sauerkraut served with bread dumplings"
$ws.Range("C3").Value = "Krůtí prsa s bramborem/ Old-czech style turkey breast with jasmine rice"
$ws.Range("C4").Value = "Koprová omáčka se sázeným vejcem / Grilled minced meat with roasted potatoes and spicy salad of roasted peppers"
$ws.Range("C5").Value = "Pasta La Vista / Tagliolini with beef tenderloin sprinkled with Grana Padano Cheese"
$ws.Range("C6").Value = "Losos na kmíně  / Baked zander with vegetables in butter served with parsley potatoes"

# Soups (column F)
$ws.Range("F2").Value = "Rajčatová polévka / Potato soup"
$ws.Range("F3").Value = "Kuřecí vývar / Lentil soup with sausages"
$ws.Range("F4").Value = "Pórková polévka / Beef consommé with meat and noodles"
$ws.Range("F5").Value = "Těstovinová polévka / Minestrone soup with pasta"
$ws.Range("F6").Value = "Hovězí vývarovka / Bank holiday. We do not serve daily menu."

# Selection moved from K3 to K4
$ws.Range("K4").Select()
